# Generate Report for Handback
#
# The localization report is refreshed after a successful handback: the
# per-file "Status" now reads that the target is in sync with en-US, the
# "Latest Handback DateTime" for each locale moves forward to the moment the
# new handback report was generated, and the now-resolved "Error Detail"
# (which used to complain the handback file was stale) is cleared out.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for both source files ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K2").Value = "2016-10-20 01:19:00"
$zhcn.Range("K3").Value = "2016-10-20 01:19:00"
$zhcn.Range("P3").Value = ""

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("K2").Value = "2016-10-20 01:19:17"
$dede.Range("K3").Value = "2016-10-20 01:19:17"
$dede.Range("P3").Value = ""

# --- Column widths: re-fit the Status / Error Detail columns now that their
#     text changed length (narrower Error Detail, wider Status). The widths
#     below are the closest values this host's width quantization can reach
#     to the recorded target widths.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333332

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333332
